$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 503
$ws.Range("J13").Value = 503
$ws.Range("L13").Value = 503
$ws.Range("N13").Value = -841
$ws.Range("H112").Value = 4388.7427
$ws.Range("J112").Value = 4595.9395
$ws.Range("L112").Value = 13787.8185
$ws.Range("N112").Value = -16003.8185
$ws.Range("H132").Value = 3925.4333
$ws.Range("I132").Value = 3391.04
$ws.Range("J132").Value = 6597.4
$ws.Range("K132").Value = 10173.12
$ws.Range("L132").Value = 19792.2
$ws.Range("M132").Value = -7643.119999999999
$ws.Range("N132").Value = -24852.2
$ws.Range("H137").Value = 8775992
$ws.Range("I137").Value = 19237118
$ws.Range("J137").Value = 2145.2903
$ws.Range("K137").Value = 57711354
$ws.Range("L137").Value = 6435.8709
$ws.Range("M137").Value = -57708804
$ws.Range("N137").Value = -11535.8709
$ws.Range("H138").Value = 3456.6038
$ws.Range("J138").Value = 3473.4443
$ws.Range("L138").Value = 10420.3329
$ws.Range("N138").Value = -20700.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6981.076
$ws.Range("I32").Value = 3687.5693
$ws.Range("J32").Value = 40857.145
$ws.Range("K32").Value = 3687.5693
$ws.Range("L32").Value = 40857.145
$ws.Range("M32").Value = -3400.5693
$ws.Range("N32").Value = -41431.145
$ws.Range("H61").Value = 2456.7441
$ws.Range("I61").Value = 2134.4856
$ws.Range("K61").Value = 2134.4856
$ws.Range("M61").Value = -1922.4856
$ws.Range("H74").Value = 6670.877
$ws.Range("I74").Value = 4527.478
$ws.Range("J74").Value = 15634.182
$ws.Range("K74").Value = 4527.478
$ws.Range("L74").Value = 15634.182
$ws.Range("M74").Value = -3653.478
$ws.Range("N74").Value = -17382.182
$ws.Range("H77").Value = 6670.877
$ws.Range("I77").Value = 4527.478
$ws.Range("J77").Value = 15634.182
$ws.Range("K77").Value = 22637.39
$ws.Range("L77").Value = 78170.91
$ws.Range("M77").Value = -18269.39
$ws.Range("N77").Value = -86906.91
$ws.Range("H102").Value = 4507.3335
$ws.Range("I102").Value = 5226.7144
$ws.Range("K102").Value = 5226.7144
$ws.Range("M102").Value = -3604.7144
$ws.Range("H136").Value = 2456.7441
$ws.Range("I136").Value = 2134.4856
$ws.Range("K136").Value = 6403.4568
$ws.Range("M136").Value = -3853.4568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2025.2069
$ws.Range("I20").Value = 1984.579
$ws.Range("K20").Value = 1984.579
$ws.Range("M20").Value = -1737.579
$ws.Range("H105").Value = 2250.5
$ws.Range("I105").Value = 2071.6667
$ws.Range("J105").Value = 2357.8
$ws.Range("K105").Value = 2071.6667
$ws.Range("L105").Value = 2357.8
$ws.Range("M105").Value = -324.6667000000002
$ws.Range("N105").Value = -5851.8
$ws.Range("H123").Value = 70780
$ws.Range("J123").Value = 70780
$ws.Range("L123").Value = 70780
$ws.Range("H134").Value = 1054.1702
$ws.Range("I134").Value = 966.4651
$ws.Range("J134").Value = 1997
$ws.Range("K134").Value = 2899.3953
$ws.Range("L134").Value = 5991
$ws.Range("M134").Value = -364.3953000000001
$ws.Range("N134").Value = -11061
$ws.Range("N123").Value = -80580

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 17020.3
$ws.Range("J28").Value = 17020.3
$ws.Range("L28").Value = 17020.3
$ws.Range("N28").Value = -17510.3
$ws.Range("H62").Value = 3826.25
$ws.Range("J62").Value = 3826.25
$ws.Range("L62").Value = 3826.25
$ws.Range("N62").Value = -5074.25
$ws.Range("H65").Value = 3826.25
$ws.Range("J65").Value = 3826.25
$ws.Range("L65").Value = 19131.25
$ws.Range("N65").Value = -25371.25
$ws.Range("H99").Value = 12702047
$ws.Range("I99").Value = 5559338
$ws.Range("J99").Value = 55558300
$ws.Range("K99").Value = 5559338
$ws.Range("L99").Value = 55558300
$ws.Range("M99").Value = -5557840
$ws.Range("N99").Value = -55561296
$ws.Range("H105").Value = 15019.9
$ws.Range("I105").Value = 17441.117
$ws.Range("K105").Value = 17441.117
$ws.Range("M105").Value = -15694.117
$ws.Range("H126").Value = 12702047
$ws.Range("I126").Value = 5559338
$ws.Range("J126").Value = 55558300
$ws.Range("K126").Value = 16678014
$ws.Range("L126").Value = 166674900
$ws.Range("M126").Value = -16675544
$ws.Range("N126").Value = -166679840
$ws.Range("H132").Value = 2272.7334
$ws.Range("I132").Value = 1892.9286
$ws.Range("K132").Value = 5678.7858
$ws.Range("M132").Value = -3148.7858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 728.9
$ws.Range("J5").Value = 644.8333
$ws.Range("L5").Value = 1934.4999
$ws.Range("N5").Value = -2158.4999
$ws.Range("H39").Value = 8777.556
$ws.Range("J39").Value = 8777.556
$ws.Range("L39").Value = 26332.668
$ws.Range("N39").Value = -26920.668
$ws.Range("H86").Value = 753.5454999999999
$ws.Range("I86").Value = 674.75
$ws.Range("J86").Value = 798.5714
$ws.Range("K86").Value = 2024.25
$ws.Range("L86").Value = 2395.7142
$ws.Range("M86").Value = -838.25
$ws.Range("N86").Value = -4767.7142
$ws.Range("H89").Value = 753.5454999999999
$ws.Range("I89").Value = 674.75
$ws.Range("J89").Value = 798.5714
$ws.Range("K89").Value = 6072.75
$ws.Range("L89").Value = 7187.1426
$ws.Range("M89").Value = -144.75
$ws.Range("N89").Value = -19043.1426
$ws.Range("H92").Value = 1834.9333
$ws.Range("I92").Value = 819.6667
$ws.Range("J92").Value = 2088.75
$ws.Range("K92").Value = 2459.0001
$ws.Range("L92").Value = 6266.25
$ws.Range("M92").Value = -1211.0001
$ws.Range("N92").Value = -8762.25
$ws.Range("H113").Value = 9806067
$ws.Range("I113").Value = 2277.923
$ws.Range("J113").Value = 15875079
$ws.Range("K113").Value = 6833.768999999999
$ws.Range("L113").Value = 47625237
$ws.Range("M113").Value = -4663.768999999999
$ws.Range("N113").Value = -47629577
$ws.Range("H122").Value = 2356.5
$ws.Range("J122").Value = 843
$ws.Range("L122").Value = 7587
$ws.Range("N122").Value = -12487
$ws.Range("H135").Value = 728.9
$ws.Range("J135").Value = 644.8333
$ws.Range("L135").Value = 5803.4997
$ws.Range("N135").Value = -10873.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5767.88
$ws.Range("I126").Value = 5994.8823
$ws.Range("J126").Value = 5285.5
$ws.Range("K126").Value = 17984.6469
$ws.Range("L126").Value = 15856.5
$ws.Range("M126").Value = -15514.6469
$ws.Range("N126").Value = -20796.5
$ws.Range("H139").Value = 97499.75
$ws.Range("J139").Value = 97499.75
$ws.Range("L139").Value = 97499.75
$ws.Range("N139").Value = -107779.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 90916264
$ws.Range("I7").Value = 111117656
$ws.Range("K7").Value = 111117656
$ws.Range("M7").Value = -111117544
$ws.Range("H46").Value = 7938302.5
$ws.Range("I46").Value = 41667390
$ws.Range("K46").Value = 41667390
$ws.Range("M46").Value = -41667202
$ws.Range("H61").Value = 5770.2856
$ws.Range("I61").Value = 6378.4
$ws.Range("J61").Value = 4250
$ws.Range("K61").Value = 6378.4
$ws.Range("L61").Value = 4250
$ws.Range("M61").Value = -6176.4
$ws.Range("H93").Value = 2579.9048
$ws.Range("I93").Value = 1466.25
$ws.Range("J93").Value = 6143.6
$ws.Range("K93").Value = 1466.25
$ws.Range("L93").Value = 6143.6
$ws.Range("M93").Value = -218.25
$ws.Range("N93").Value = -8639.6
$ws.Range("H98").Value = 20000
$ws.Range("J98").Value = 20000
$ws.Range("L98").Value = 20000
$ws.Range("N98").Value = -25990
$ws.Range("H113").Value = 5770.2856
$ws.Range("I113").Value = 6378.4
$ws.Range("J113").Value = 4250
$ws.Range("K113").Value = 6378.4
$ws.Range("L113").Value = 4250
$ws.Range("M113").Value = -4208.4
$ws.Range("H122").Value = 7032.654
$ws.Range("I122").Value = 6419.65
$ws.Range("J122").Value = 9076
$ws.Range("K122").Value = 19258.95
$ws.Range("L122").Value = 27228
$ws.Range("M122").Value = -16808.95
$ws.Range("N122").Value = -32128
$ws.Range("H126").Value = 90916264
$ws.Range("I126").Value = 111117656
$ws.Range("K126").Value = 333352968
$ws.Range("M126").Value = -333350498
$ws.Range("N61").Value = -4654
$ws.Range("N113").Value = -8590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 31051.666
$ws.Range("I62").Value = 33747.5
$ws.Range("J62").Value = 17572.5
$ws.Range("K62").Value = 33747.5
$ws.Range("L62").Value = 17572.5
$ws.Range("M62").Value = -33123.5
$ws.Range("N62").Value = -18820.5
$ws.Range("H65").Value = 31051.666
$ws.Range("I65").Value = 33747.5
$ws.Range("J65").Value = 17572.5
$ws.Range("K65").Value = 168737.5
$ws.Range("L65").Value = 87862.5
$ws.Range("M65").Value = -165617.5
$ws.Range("N65").Value = -94102.5
$ws.Range("H113").Value = 429.48148
$ws.Range("I113").Value = 376.86365
$ws.Range("K113").Value = 1130.59095
$ws.Range("M113").Value = 1039.40905
$ws.Range("H126").Value = 111114510
$ws.Range("I126").Value = 41670404
$ws.Range("K126").Value = 125011212
$ws.Range("M126").Value = -125008742
$ws.Range("H132").Value = 2729.543
$ws.Range("I132").Value = 2577.2415
$ws.Range("J132").Value = 3465.6667
$ws.Range("K132").Value = 7731.7245
$ws.Range("L132").Value = 10397.0001
$ws.Range("M132").Value = -5201.7245
$ws.Range("N132").Value = -15457.0001
